$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ProjectID"
$ws.Range("C1").Value = "ProjectDependency"

$ws.Range("C1").Select()
